$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-01-27 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-01-28 Tuesday", 2)

$d.Content.Find.Execute("29×55=", $true, $false, $false, $false, $false, $true, 1, $false, "73×24=", 2)
$d.Content.Find.Execute("17×69=", $true, $false, $false, $false, $false, $true, 1, $false, "67×21=", 2)
$d.Content.Find.Execute("16×91=", $true, $false, $false, $false, $false, $true, 1, $false, "84×21=", 2)
$d.Content.Find.Execute("63×23=", $true, $false, $false, $false, $false, $true, 1, $false, "20×70=", 2)
$d.Content.Find.Execute("31×92=", $true, $false, $false, $false, $false, $true, 1, $false, "91×33=", 2)

$d.Content.Find.Execute("57×30=", $true, $false, $false, $false, $false, $true, 1, $false, "94×13=", 2)
$d.Content.Find.Execute("72×36=", $true, $false, $false, $false, $false, $true, 1, $false, "83×14=", 2)
$d.Content.Find.Execute("46×84=", $true, $false, $false, $false, $false, $true, 1, $false, "44×53=", 2)
$d.Content.Find.Execute("24×84=", $true, $false, $false, $false, $false, $true, 1, $false, "82×82=", 2)
$d.Content.Find.Execute("37×49=", $true, $false, $false, $false, $false, $true, 1, $false, "70×83=", 2)

$d.Content.Find.Execute("99×21=", $true, $false, $false, $false, $false, $true, 1, $false, "92×68=", 2)
$d.Content.Find.Execute("25×58=", $true, $false, $false, $false, $false, $true, 1, $false, "61×15=", 2)
$d.Content.Find.Execute("46×28=", $true, $false, $false, $false, $false, $true, 1, $false, "11×64=", 2)
$d.Content.Find.Execute("21×17=", $true, $false, $false, $false, $false, $true, 1, $false, "37×65=", 2)
$d.Content.Find.Execute("76×39=", $true, $false, $false, $false, $false, $true, 1, $false, "52×53=", 2)

$d.Content.Find.Execute("99×83=", $true, $false, $false, $false, $false, $true, 1, $false, "40×51=", 2)
$d.Content.Find.Execute("65×89=", $true, $false, $false, $false, $false, $true, 1, $false, "98×90=", 2)
$d.Content.Find.Execute("50×79=", $true, $false, $false, $false, $false, $true, 1, $false, "85×15=", 2)
$d.Content.Find.Execute("89×92=", $true, $false, $false, $false, $false, $true, 1, $false, "87×11=", 2)
$d.Content.Find.Execute("87×17=", $true, $false, $false, $false, $false, $true, 1, $false, "91×42=", 2)

$d.Content.Find.Execute("77×32=", $true, $false, $false, $false, $false, $true, 1, $false, "96×40=", 2)
$d.Content.Find.Execute("80×98=", $true, $false, $false, $false, $false, $true, 1, $false, "65×60=", 2)
$d.Content.Find.Execute("29×93=", $true, $false, $false, $false, $false, $true, 1, $false, "15×14=", 2)
$d.Content.Find.Execute("85×47=", $true, $false, $false, $false, $false, $true, 1, $false, "54×80=", 2)
$d.Content.Find.Execute("53×64=", $true, $false, $false, $false, $false, $true, 1, $false, "95×68=", 2)
